$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 5")
$r = $ws.Range("A20")
Write-Output ($r | Get-Member | Out-String)
